$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 1005
$ws.Range("I3").Value = 1005
$ws.Range("K3").Value = 1005
$ws.Range("M3").Value = -890
# Row 21
$ws.Range("H21").Value = 80000
$ws.Range("J21").Value = 80000
$ws.Range("L21").Value = 80000
$ws.Range("N21").Value = -80936
# Row 23
$ws.Range("H23").Value = 80000
$ws.Range("J23").Value = 80000
$ws.Range("L23").Value = 80000
$ws.Range("N23").Value = -80468
# Row 29
$ws.Range("H29").Value = 4000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 32
$ws.Range("H32").Value = 864.8182
$ws.Range("I32").Value = 727.75
$ws.Range("J32").Value = 943.1429000000001
$ws.Range("K32").Value = 727.75
$ws.Range("L32").Value = 943.1429000000001
$ws.Range("M32").Value = -401.75
$ws.Range("N32").Value = -1595.1429
# Row 33
$ws.Range("H33").Value = 532.96155
$ws.Range("I33").Value = 172.53334
$ws.Range("J33").Value = 1024.4546
$ws.Range("K33").Value = 172.53334
$ws.Range("L33").Value = 1024.4546
$ws.Range("M33").Value = 56.46665999999999
$ws.Range("N33").Value = -1482.4546
# Row 38
$ws.Range("H38").Value = 476.66666
$ws.Range("I38").Value = 145.45454
$ws.Range("J38").Value = 997.1429000000001
$ws.Range("K38").Value = 436.36362
$ws.Range("L38").Value = 2991.4287
$ws.Range("M38").Value = -64.36362000000003
$ws.Range("N38").Value = -3735.4287
# Row 58
$ws.Range("H58").Value = 996.6667
$ws.Range("I58").Value = 996.6667
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2990.0001
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2840.0001
$ws.Range("N58").ClearContents()
# Row 87
$ws.Range("H87").Value = 33966.57
$ws.Range("J87").Value = 33966.57
$ws.Range("L87").Value = 33966.57
$ws.Range("N87").Value = -36462.57
# Row 90
$ws.Range("H90").Value = 33966.57
$ws.Range("J90").Value = 33966.57
$ws.Range("L90").Value = 101899.71
$ws.Range("N90").Value = -114379.71
# Row 131
$ws.Range("H131").Value = 1133.1818
$ws.Range("I131").Value = 637
$ws.Range("J131").Value = 1546.6666
$ws.Range("K131").Value = 1911
$ws.Range("L131").Value = 4639.9998
$ws.Range("M131").Value = 3129
$ws.Range("N131").Value = -14719.9998
# Row 132
$ws.Range("H132").Value = 744033.8
$ws.Range("I132").Value = 1376.8214
$ws.Range("J132").Value = 4902913
$ws.Range("K132").Value = 4130.4642
$ws.Range("L132").Value = 14708739
$ws.Range("M132").Value = -1600.4642
$ws.Range("N132").Value = -14713799
# Row 135
$ws.Range("H135").Value = 20942.666
$ws.Range("I135").Value = 24052.744
$ws.Range("J135").Value = 4226
$ws.Range("K135").Value = 216474.696
$ws.Range("L135").Value = 38034
$ws.Range("M135").Value = -213939.696
$ws.Range("N135").Value = -43104
# Row 137
$ws.Range("H137").Value = 2274255.2
$ws.Range("I137").Value = 3847380
$ws.Range("J137").Value = 1963.8889
$ws.Range("K137").Value = 11542140
$ws.Range("L137").Value = 5891.6667
$ws.Range("M137").Value = -11539590
$ws.Range("N137").Value = -10991.6667
# Row 138
$ws.Range("H138").Value = 2826516
$ws.Range("I138").Value = 1185.2703
$ws.Range("J138").Value = 7578208.5
$ws.Range("K138").Value = 3555.810899999999
$ws.Range("L138").Value = 22734625.5
$ws.Range("M138").Value = 1584.189100000001
$ws.Range("N138").Value = -22744905.5
# Row 141
$ws.Range("H141").Value = 1907.6774
$ws.Range("I141").Value = 1921.2667
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 5763.800099999999
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = -583.8000999999995
$ws.Range("N141").Value = -14860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 7733.1333
$ws.Range("I26").Value = 6399.4
$ws.Range("J26").Value = 8400
$ws.Range("K26").Value = 6399.4
$ws.Range("L26").Value = 8400
$ws.Range("M26").Value = -6069.4
$ws.Range("N26").Value = -9060
# Row 32
$ws.Range("H32").Value = 1493.15
$ws.Range("I32").Value = 1005.13336
$ws.Range("J32").Value = 2957.2
$ws.Range("K32").Value = 1005.13336
$ws.Range("L32").Value = 2957.2
$ws.Range("M32").Value = -718.13336
$ws.Range("N32").Value = -3531.2
# Row 35
$ws.Range("H35").Value = 5250.3
$ws.Range("I35").Value = 1750
$ws.Range("J35").Value = 6125.375
$ws.Range("K35").Value = 1750
$ws.Range("L35").Value = 6125.375
$ws.Range("M35").Value = -1344
$ws.Range("N35").Value = -6937.375
# Row 124
$ws.Range("H124").Value = 30429
$ws.Range("J124").Value = 30429
$ws.Range("L124").Value = 30429
$ws.Range("N124").Value = -40249
# Row 125
$ws.Range("H125").Value = 55000
$ws.Range("J125").Value = 55000
$ws.Range("L125").Value = 55000
$ws.Range("N125").Value = -64840
# Row 132
$ws.Range("H132").Value = 62193.41
$ws.Range("I132").Value = 35071.1
$ws.Range("J132").Value = 219502.8
$ws.Range("K132").Value = 105213.3
$ws.Range("L132").Value = 658508.3999999999
$ws.Range("M132").Value = -102683.3
$ws.Range("N132").Value = -663568.3999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 124
$ws.Range("H124").Value = 34885
$ws.Range("J124").Value = 34885
$ws.Range("L124").Value = 34885
$ws.Range("N124").Value = -44705
# Row 134
$ws.Range("H134").Value = 1506.7413
$ws.Range("I134").Value = 731.3953
$ws.Range("J134").Value = 3729.4
$ws.Range("K134").Value = 2194.1859
$ws.Range("L134").Value = 11188.2
$ws.Range("M134").Value = 340.8141000000001
$ws.Range("N134").Value = -16258.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3130.2888
$ws.Range("I31").Value = 1292.0385
$ws.Range("J31").Value = 5645.7896
$ws.Range("K31").Value = 1292.0385
$ws.Range("L31").Value = 5645.7896
$ws.Range("M31").Value = -997.0385000000001
$ws.Range("N31").Value = -6235.7896
# Row 34
$ws.Range("H34").Value = 3130.2888
$ws.Range("I34").Value = 1292.0385
$ws.Range("J34").Value = 5645.7896
$ws.Range("K34").Value = 1292.0385
$ws.Range("L34").Value = 5645.7896
$ws.Range("M34").Value = -1090.0385
$ws.Range("N34").Value = -6049.7896
# Row 58
$ws.Range("H58").Value = 33335740
$ws.Range("I58").Value = 41668820
$ws.Range("J58").Value = 3418.8333
$ws.Range("K58").Value = 41668820
$ws.Range("L58").Value = 3418.8333
$ws.Range("M58").Value = -41668617
$ws.Range("N58").Value = -3824.8333
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 124
$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -24910
# Row 132
$ws.Range("H132").Value = 40904.52
$ws.Range("I132").Value = 22113.318
$ws.Range("J132").Value = 335300
$ws.Range("K132").Value = 66339.954
$ws.Range("L132").Value = 1005900
$ws.Range("M132").Value = -63809.954
$ws.Range("N132").Value = -1010960
# Row 134
$ws.Range("H134").Value = 23993.229
$ws.Range("I134").Value = 1516.4286
$ws.Range("K134").Value = 4549.2858
$ws.Range("M134").Value = -2014.2858
# Row 136
$ws.Range("H136").Value = 33335740
$ws.Range("I136").Value = 41668820
$ws.Range("J136").Value = 3418.8333
$ws.Range("K136").Value = 125006460
$ws.Range("L136").Value = 10256.4999
$ws.Range("M136").Value = -125003910
$ws.Range("N136").Value = -15356.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2701.5789
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2701.5789
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 8104.736699999999
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -8328.736699999999
# Row 113
$ws.Range("H113").Value = 514.8936
$ws.Range("J113").Value = 763.63635
$ws.Range("L113").Value = 2290.90905
$ws.Range("N113").Value = -6630.90905
# Row 131
$ws.Range("H131").Value = 894.2308
$ws.Range("I131").Value = 435.5
$ws.Range("J131").Value = 1063.2368
$ws.Range("K131").Value = 1306.5
$ws.Range("L131").Value = 3189.7104
$ws.Range("M131").Value = 3733.5
$ws.Range("N131").Value = -13269.7104

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2117.8333
$ws.Range("I126").Value = 1466.6666
$ws.Range("J126").Value = 2769
$ws.Range("K126").Value = 4399.9998
$ws.Range("L126").Value = 8307
$ws.Range("M126").Value = -1929.9998
$ws.Range("N126").Value = -13247

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1835.875
$ws.Range("J7").Value = 2260
$ws.Range("L7").Value = 2260
$ws.Range("N7").Value = -2484
# Row 126
$ws.Range("H126").Value = 1835.875
$ws.Range("J126").Value = 2260
$ws.Range("L126").Value = 6780
$ws.Range("N126").Value = -11720
# Row 127
$ws.Range("H127").Value = 47500
$ws.Range("J127").Value = 47500
$ws.Range("L127").Value = 47500
$ws.Range("N127").Value = -57420
# Row 132
$ws.Range("H132").Value = 22806.986
$ws.Range("I132").Value = 9314.234
$ws.Range("J132").Value = 146169.28
$ws.Range("K132").Value = 27942.702
$ws.Range("L132").Value = 438507.84
$ws.Range("M132").Value = -25412.702
$ws.Range("N132").Value = -443567.84

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 132
$ws.Range("H132").Value = 54212.484
$ws.Range("I132").Value = 38485.168
$ws.Range("J132").Value = 220922
$ws.Range("K132").Value = 115455.504
$ws.Range("L132").Value = 662766
$ws.Range("M132").Value = -112925.504
$ws.Range("N132").Value = -667826
